$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old rows 341:344 (previous week's data) get pushed down to become
# rows 345:348, preserving their original values/styles. A new week's
# data (rows 341:344) is written in their place.
$ws.Range("A341:R344").Copy($ws.Range("A345"))

# --- Row 341: updated with new week's data ---
$ws.Range("D341").Value = 44595
$ws.Range("J341").Value = 800
$ws.Range("K341").Value = 16000
$ws.Range("L341").Value = 17000
$ws.Range("M341").Value = 16500
$ws.Range("O341").Value = "Región de O'Higgins"
$ws.Range("P341").Value = 917

# --- Row 342: updated with new week's data ---
$ws.Range("D342").Value = 44595
$ws.Range("J342").Value = 400
$ws.Range("K342").Value = 15000
$ws.Range("L342").Value = 15000
$ws.Range("M342").Value = 15000
$ws.Range("O342").Value = "Región de O'Higgins"
$ws.Range("P342").Value = 833

# --- Row 343: updated with new week's data ---
$ws.Range("D343").Value = 44595
$ws.Range("H343").Value = "Semiduro"
$ws.Range("L343").Value = 10000
$ws.Range("M343").Value = 9500
$ws.Range("N343").Value = "$/bandeja 18 kilos"
$ws.Range("O343").Value = "Región de O'Higgins"
$ws.Range("P343").Value = 528
$ws.Range("Q343").Value = 18

# --- Row 344: updated with new week's data ---
$ws.Range("D344").Value = 44595
$ws.Range("H344").Value = "Semiduro"
$ws.Range("N344").Value = "$/bandeja 18 kilos"
$ws.Range("O344").Value = "Región de O'Higgins"
$ws.Range("P344").Value = 444
$ws.Range("Q344").Value = 18
